$wb = $excel.ActiveWorkbook

# Row 41 (G=5478) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1322.7646
$ws.Range("I41").Value = 2151.5
$ws.Range("K41").Value = 2151.5
$ws.Range("M41").Value = -1711.5

# Row 53 (G=5479) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1232.2
$ws.Range("I53").Value = 2207.4
$ws.Range("J53").Value = 257
$ws.Range("K53").Value = 2207.4
$ws.Range("L53").Value = 257
$ws.Range("M53").Value = -1570.4
$ws.Range("N53").Value = -1531

# Row 80 (G=12605) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2140.4
$ws.Range("I80").Value = 2743.4285
$ws.Range("J80").Value = 733.3333
$ws.Range("K80").Value = 8230.2855
$ws.Range("L80").Value = 2199.9999
$ws.Range("M80").Value = -7232.2855
$ws.Range("N80").Value = -4195.9999

# Row 83 (G=12605) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2140.4
$ws.Range("I83").Value = 2743.4285
$ws.Range("J83").Value = 733.3333
$ws.Range("K83").Value = 24690.8565
$ws.Range("L83").Value = 6599.9997
$ws.Range("M83").Value = -19698.8565
$ws.Range("N83").Value = -16583.9997

# Row 129 (G=36115) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 900.6269
$ws.Range("J129").Value = 882.18463
$ws.Range("L129").Value = 2646.55389
$ws.Range("N129").Value = -12646.55389

# Row 137 (G=44013) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1401.76
$ws.Range("I137").Value = 1320.1818
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 3960.5454
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -1410.5454
$ws.Range("N137").Value = -11100

# Row 138 (G=44169) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3021.1714
$ws.Range("I138").Value = 2675.5
$ws.Range("J138").Value = 4019.7778
$ws.Range("K138").Value = 8026.5
$ws.Range("L138").Value = 12059.3334
$ws.Range("M138").Value = -2886.5
$ws.Range("N138").Value = -22339.3334

# Row 141 (G=44161) on sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1402775.8
$ws.Range("I141").Value = 1868502.4
$ws.Range("K141").Value = 5605507.199999999
$ws.Range("M141").Value = -5600327.199999999

# Row 32 (G=44147) on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3621.3274
$ws.Range("I32").Value = 2661.3777
$ws.Range("J32").Value = 7941.1
$ws.Range("K32").Value = 2661.3777
$ws.Range("L32").Value = 7941.1
$ws.Range("M32").Value = -2374.3777
$ws.Range("N32").Value = -8515.1

# Row 102 (G=19945) on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1763.5385
$ws.Range("I102").Value = 1618.8334
$ws.Range("K102").Value = 1618.8334
$ws.Range("M102").Value = 3.166600000000017

# Row 132 (G=43997) on sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1820.35
$ws.Range("I132").Value = 1091.25
$ws.Range("K132").Value = 3273.75
$ws.Range("M132").Value = -743.75

# Row 107 (G=27706) on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2953.125
$ws.Range("I107").Value = 3520
$ws.Range("K107").Value = 3520
$ws.Range("M107").Value = -1600

# Row 134 (G=43998) on sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14741.762
$ws.Range("I134").Value = 14681.059
$ws.Range("J134").Value = 14999.75
$ws.Range("K134").Value = 44043.177
$ws.Range("L134").Value = 44999.25
$ws.Range("M134").Value = -41508.177
$ws.Range("N134").Value = -50069.25

# Row 31 (G=44023) on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3762
$ws.Range("I31").Value = 3003.3333
$ws.Range("K31").Value = 3003.3333
$ws.Range("M31").Value = -2708.3333

# Row 34 (G=44023) on sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3762
$ws.Range("I34").Value = 3003.3333
$ws.Range("K34").Value = 3003.3333
$ws.Range("M34").Value = -2801.3333

# Row 37 (G=9516) on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 100000
$ws.Range("J37").Value = 100000
$ws.Range("L37").Value = 300000
$ws.Range("N37").Value = -300224

# Row 131 (G=36060) on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 756.62
$ws.Range("J131").Value = 783.16486
$ws.Range("L131").Value = 2349.49458
$ws.Range("N131").Value = -12429.49458

# Row 137 (G=44088) on sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3139.2173
$ws.Range("J137").Value = 3162.2104
$ws.Range("L137").Value = 9486.6312
$ws.Range("N137").Value = -19686.6312

# Row 70 (G=14146) on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8752.647000000001
$ws.Range("I70").Value = 11880
$ws.Range("J70").Value = 4285
$ws.Range("K70").Value = 11880
$ws.Range("L70").Value = 4285
$ws.Range("M70").Value = -11610
$ws.Range("N70").Value = -4825

# Row 73 (G=14146) on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8752.647000000001
$ws.Range("I73").Value = 11880
$ws.Range("J73").Value = 4285
$ws.Range("K73").Value = 11880
$ws.Range("L73").Value = 4285
$ws.Range("M73").Value = -10944
$ws.Range("N73").Value = -6157

# Row 80 (G=12521) on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3014.9092
$ws.Range("I80").Value = 2931.6667
$ws.Range("J80").Value = 3114.8
$ws.Range("K80").Value = 2931.6667
$ws.Range("L80").Value = 3114.8
$ws.Range("M80").Value = -1933.6667
$ws.Range("N80").Value = -5110.8

# Row 83 (G=12521) on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3014.9092
$ws.Range("I83").Value = 2931.6667
$ws.Range("J83").Value = 3114.8
$ws.Range("K83").Value = 14658.3335
$ws.Range("L83").Value = 15574
$ws.Range("M83").Value = -9666.333500000001
$ws.Range("N83").Value = -25558

# Row 122 (G=36182) on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2502.3333
$ws.Range("J122").Value = 2644
$ws.Range("L122").Value = 7932
$ws.Range("N122").Value = -12832

# Row 136 (G=42218) on sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 9865.25
$ws.Range("J136").Value = 9865.25
$ws.Range("L136").Value = 29595.75
$ws.Range("N136").Value = -34695.75

# Row 22 (G=5277) on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2737.25
$ws.Range("I22").Value = 5300
$ws.Range("J22").Value = 1883
$ws.Range("K22").Value = 5300
$ws.Range("L22").Value = 1883
$ws.Range("M22").Value = -5005
$ws.Range("N22").Value = -2473

# Row 27 (G=5277) on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2737.25
$ws.Range("I27").Value = 5300
$ws.Range("J27").Value = 1883
$ws.Range("K27").Value = 5300
$ws.Range("L27").Value = 1883
$ws.Range("M27").Value = -5193
$ws.Range("N27").Value = -2097

# Row 68 (G=12563) on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2514.5386
$ws.Range("I68").Value = 2244.4546
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 2244.4546
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -1495.4546
$ws.Range("N68").Value = -5498

# Row 71 (G=12563) on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2514.5386
$ws.Range("I71").Value = 2244.4546
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 11222.273
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -7478.273000000001
$ws.Range("N71").Value = -27488

# Row 100 (G=19995) on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1328.3334
$ws.Range("I100").Value = 1297.5
$ws.Range("J100").Value = 1390
$ws.Range("K100").Value = 1297.5
$ws.Range("L100").Value = 1390
$ws.Range("M100").Value = -756.5
$ws.Range("N100").Value = -2472

# Row 132 (G=44058) on sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2160.4614
$ws.Range("I132").Value = 2220.9092
$ws.Range("K132").Value = 6662.7276
$ws.Range("M132").Value = -4132.7276

# Row 7 (G=2661) on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 43504.5
$ws.Range("I7").Value = 7004
$ws.Range("K7").Value = 7004
$ws.Range("M7").Value = -6891

# Row 119 (G=26289) on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# Row 120 (G=26310) on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

# Row 132 (G=44029) on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1843
$ws.Range("I132").Value = 1330.4
$ws.Range("K132").Value = 3991.2
$ws.Range("M132").Value = -1461.2

# Row 139 (G=43312) on sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 58887.777
$ws.Range("J139").Value = 59998.75
$ws.Range("L139").Value = 59998.75
$ws.Range("N139").Value = -70278.75
